$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.109.63'
$ws.Range('E2').Value = '  +2.24%  '
$ws.Range('D3').Value = '3.459.35'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.14'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.60'
$ws.Range('E6').Value = '  +3.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.609'
$ws.Range('E7').Value = '  +6.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.467.93'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.27'
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.452'
$ws.Range('E12').Value = '  +2.46%  '
$ws.Range('D13').Value = '4.057.51'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('E14').Value = '  +0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000191'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '28.24'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '65.134.32'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').Value = '3.473.81'
$ws.Range('E18').Value = '  +2.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.47'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.33'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.98'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.17'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.555'
$ws.Range('E23').Value = '  +3.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.26'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.08'
$ws.Range('E27').Value = '  +5.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('E30').Value = '  +9.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.17'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.66'
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.29'
$ws.Range('E34').Value = '  +6.01%  '
$ws.Range('E35').Value = '  +11.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.67'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('E37').Value = '  +5.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0783'
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.86'
$ws.Range('E39').Value = '  +6.88%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.908.31'
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.04'
$ws.Range('E41').Value = '  +1.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.71'
$ws.Range('E42').Value = '  +7.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.18'
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0318'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.781'
$ws.Range('E45').Value = '  +3.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.95'
$ws.Range('E46').Value = '  +11.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '324.52'
$ws.Range('E47').Value = '  +11.58%  '
$ws.Range('E48').Value = '  +3.04%  '
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.874'
$ws.Range('E50').Value = '  +4.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.18'
$ws.Range('E51').Value = '  -0.27%  '
